$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C: "Data Tiparire" (C) shifts to D,
# "Perioada Internarii" (D) shifts to E, "Urgenta" (E) shifts to F.
$ws.Columns.Item(3).Insert()

# New "Varsta" header in the freshly inserted column C.
$ws.Cells.Item(1, 3).Value = "Varsta"

# New "LDL COLESTEROL" column header in G, copying A1's header formatting
# (bold, centered, bordered) onto it first so the style index is reused
# instead of creating a brand-new one.
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 7))
$ws.Cells.Item(1, 7).Value = "LDL COLESTEROL"

# Fill in the age values for the 5 data rows.
$ws.Cells.Item(2, 3).Value = "71 ani si 0 luni `n"
$ws.Cells.Item(3, 3).Value = "78 ani si 2 luni `n"
$ws.Cells.Item(4, 3).Value = "55 ani si 3 luni `n"
$ws.Cells.Item(5, 3).Value = "67 ani si 1 luni `n"
$ws.Cells.Item(6, 3).Value = "64 ani si 6 luni `n"

# Writing multi-line text auto-expands the row height; AutoFit puts each
# row back to the sheet's default (and clears the "custom height" flag).
$ws.Range("2:6").AutoFit()

# Extend the used range down to row 51 / out to column G: register the
# (empty) "LDL COLESTEROL" cells for the existing 5 data rows, and the
# blank trailing rows 7-51 across all 7 columns, by tiling a single blank
# cell (taken from well outside the used range) over the destinations.
$blank = $ws.Cells.Item(100, 1)
$blank.Copy($ws.Range("G2:G6"))
$blank.Copy($ws.Range("A7:G51"))

$ws.Range("A1").Select() | Out-Null
